$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the transformation parameters (Hshift, Hstrech, Vstrech, Vshift)
$ws.Range("R27").Value = 3
$ws.Range("R29").Value = 4
$ws.Range("R33").Value = 3
$ws.Range("R37").Value = 5

# Move the active selection to match the saved view state
$ws.Range("R38").Select()
